# Replace the Kahraman-model Solar notified-production series (2025-08-?? / day 45747)
# with the new model's output, shifted four days forward (day 45751..45752).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column-A timestamps (Excel serial date/time, 15-min resolution)
$aValues = @(
    45751.01041666666, 45751.02083333334, 45751.03125, 45751.04166666666, 45751.05208333334, 45751.0625,
    45751.07291666666, 45751.08333333334, 45751.09375, 45751.10416666666, 45751.11458333334, 45751.125,
    45751.13541666666, 45751.14583333334, 45751.15625, 45751.16666666666, 45751.17708333334, 45751.1875,
    45751.19791666666, 45751.20833333334, 45751.21875, 45751.22916666666, 45751.23958333334, 45751.25,
    45751.26041666666, 45751.27083333334, 45751.28125, 45751.29166666666, 45751.30208333334, 45751.3125,
    45751.32291666666, 45751.33333333334, 45751.34375, 45751.35416666666, 45751.36458333334, 45751.375,
    45751.38541666666, 45751.39583333334, 45751.40625, 45751.41666666666, 45751.42708333334, 45751.4375,
    45751.44791666666, 45751.45833333334, 45751.46875, 45751.47916666666, 45751.48958333334, 45751.5,
    45751.51041666666, 45751.52083333334, 45751.53125, 45751.54166666666, 45751.55208333334, 45751.5625,
    45751.57291666666, 45751.58333333334, 45751.59375, 45751.60416666666, 45751.61458333334, 45751.625,
    45751.63541666666, 45751.64583333334, 45751.65625, 45751.66666666666, 45751.67708333334, 45751.6875,
    45751.69791666666, 45751.70833333334, 45751.71875, 45751.72916666666, 45751.73958333334, 45751.75,
    45751.76041666666, 45751.77083333334, 45751.78125, 45751.79166666666, 45751.80208333334, 45751.8125,
    45751.82291666666, 45751.83333333334, 45751.84375, 45751.85416666666, 45751.86458333334, 45751.875,
    45751.88541666666, 45751.89583333334, 45751.90625, 45751.91666666666, 45751.92708333334, 45751.9375,
    45751.94791666666, 45751.95833333334, 45751.96875, 45751.97916666666, 45751.98958333334, 45752
)

# New column-B notified-production values (MW)
$bValues = @(
    1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 10, 2, 3, 4, 32, 31, 38, 49, 198, 220, 256,
    305, 514, 554, 592, 638, 821, 847, 874, 898, 1087, 1102, 1115, 1126, 1245, 1252, 1254, 1253, 1240, 1235,
    1228, 1220, 1177, 1164, 1152, 1141, 1058, 1042, 1027, 1013, 891, 871, 852, 831, 606, 583, 562, 537, 272,
    249, 228, 210, 32, 23, 20, 18, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1
)

# Data starts on row 2 (row 1 is the header: Timestamp / Notified Production (MW))
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
